# DOMA-3100 add formatter convert to number for some colomns
#
# For the "processing", "completed", "canceled", "deferred", "closed" and
# "new_or_reopened" template columns (row i and row i+1 in the ticket
# analytics export) append the ":formatN()" formatter to the placeholder
# text and switch the cell's number format from text ("@") to a plain
# integer number format ("0") so the values are rendered as numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$columns = @("C", "D", "E", "F", "G", "H")
$rows = @(2, 3)

foreach ($row in $rows) {
    foreach ($col in $columns) {
        $cell = $ws.Range("$col$row")
        $oldText = $cell.Text
        $newText = $oldText -replace '\}$', ':formatN()}'
        $cell.Value = $newText
        $cell.NumberFormat = "0"
    }
}
